$d = $word.ActiveDocument
$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)
$last.Range.Text = "Esses programas precisam estar instalados o quanto antes (antes das aulas)."
